$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.175.67'
$ws.Range("E2").Value = '  +1.95%  '

# Row 3
$ws.Range("D3").Value = '2.381.83'
$ws.Range("E3").Value = '  +4.11%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''302.98'
$ws.Range("E5").Value = '  +0.76%  '

# Row 6
$ws.Range("D6").Value = '''97.03'
$ws.Range("E6").Value = '  +2.32%  '

# Row 7
$ws.Range("E7").Value = '  +0.83%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +2.13%  '

# Row 10
$ws.Range("D10").Value = '''34.28'
$ws.Range("E10").Value = '  +0.05%  '

# Row 12
$ws.Range("E12").Value = '  +2.09%  '

# Row 13
$ws.Range("D13").Value = '''18.43'
$ws.Range("E13").Value = '  -2.84%  '

# Row 14
$ws.Range("E14").Value = '  +1.44%  '

# Row 15
$ws.Range("D15").Value = '2.754.22'
$ws.Range("E15").Value = '  +4.10%  '

# Row 16
$ws.Range("D16").Value = '2.379.58'
$ws.Range("E16").Value = '  +3.63%  '

# Row 17
$ws.Range("D17").Value = '''0.808'
$ws.Range("E17").Value = '  +3.91%  '

# Row 18
$ws.Range("D18").Value = '43.184.38'
$ws.Range("E18").Value = '  +2.07%  '

# Row 19
$ws.Range("E19").Value = '  +0.26%  '

# Row 20
$ws.Range("D20").Value = '''6.32'
$ws.Range("E20").Value = '  +6.17%  '

# Row 21
$ws.Range("E21").Value = '  +0.38%  '

# Row 22
$ws.Range("D22").Value = '''68.70'
$ws.Range("E22").Value = '  +1.95%  '

# Row 23
$ws.Range("D23").Value = '''235.24'
$ws.Range("E23").Value = '  -0.08%  '

# Row 24
$ws.Range("E24").Value = '  -0.97%  '

# Row 25
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("D26").Value = '''2.44'
$ws.Range("E26").Value = '  +1.70%  '

# Row 27
$ws.Range("D27").Value = '''24.92'
$ws.Range("E27").Value = '  +3.04%  '

# Row 28
$ws.Range("E28").Value = '  +0.45%  '

# Row 29
$ws.Range("D29").Value = '''9.14'
$ws.Range("E29").Value = '  +1.32%  '

# Row 30
$ws.Range("D30").Value = '''31.52'
$ws.Range("E30").Value = '  -0.33%  '

# Row 31
$ws.Range("E31").Value = '  +0.08%  '

# Row 32
$ws.Range("E32").Value = '  +2.74%  '

# Row 33
$ws.Range("D33").Value = '''0.0736'
$ws.Range("E33").Value = '  +6.27%  '

# Row 34
$ws.Range("D34").Value = '''17.12'
$ws.Range("E34").Value = '  -2.00%  '

# Row 35
$ws.Range("E35").Value = '  +7.17%  '

# Row 36
$ws.Range("D36").Value = '''0.103'
$ws.Range("E36").Value = '  +3.24%  '

# Row 37
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '''2.30'
$ws.Range("E37").Value = '  -1.02%  '

# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''4.29'
$ws.Range("E38").Value = '  -0.95%  '

# Row 39
$ws.Range("D39").Value = '''2.80'
$ws.Range("E39").Value = '  +4.85%  '

# Row 40
$ws.Range("D40").Value = '''22.33'
$ws.Range("E40").Value = '  +12.39%  '

# Row 41
$ws.Range("E41").Value = '  +0.50%  '

# Row 42
$ws.Range("D42").Value = '''104.19'
$ws.Range("E42").Value = '  -36.67%  '

# Row 43
$ws.Range("D43").Value = '1.957.56'
$ws.Range("E43").Value = '  +0.60%  '

# Row 44
$ws.Range("E44").Value = '  +0.92%  '

# Row 45
$ws.Range("E45").Value = '  +1.83%  '

# Row 46
$ws.Range("E46").Value = '  +1.17%  '

# Row 47
$ws.Range("D47").Value = '''9.27'
$ws.Range("E47").Value = '  -9.94%  '

# Row 48
$ws.Range("D48").Value = '''52.81'
$ws.Range("E48").Value = '  +0.16%  '

# Row 49
$ws.Range("E49").Value = '  +3.48%  '

# Row 50
$ws.Range("D50").Value = '''71.88'
$ws.Range("E50").Value = '  +2.01%  '

# Row 51
$ws.Range("E51").Value = '  +1.60%  '
